# Generate Report for Handoff
# Updates the "fd424155-71a7-4b75-ae29-702ad9e1908a" file's row across all
# sheets to reflect that the file is now ready for handoff: status changes
# from "In Translation" to "Ready for handoff", Priority changes from "ht"
# to "mt", and the Latest Handoff Datetime / Latest HO Xliff Generate Date
# are refreshed with the new handoff timestamp.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-31 20:17:21"
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-31 20:17:16"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-31 20:17:21"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333
